# Apply updated Price (D) and Volume(1h) (E) values to the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.253.30'
$ws.Range('E2').Value = '  -0.92%  '
$ws.Range('D3').Value = '2.270.24'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'" + '110.73'
$ws.Range('E5').Value = '  -2.97%  '
$ws.Range('D6').Value = "'" + '263.11'
$ws.Range('E6').Value = '  -2.00%  '
$ws.Range('D7').Value = "'" + '0.646'
$ws.Range('E7').Value = '  +3.29%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('D9').Value = "'" + '0.604'
$ws.Range('E9').Value = '  -3.16%  '
$ws.Range('D10').Value = "'" + '46.29'
$ws.Range('E10').Value = '  -4.11%  '
$ws.Range('D11').Value = "'" + '0.0932'
$ws.Range('E11').Value = '  -1.84%  '
$ws.Range('D12').Value = "'" + '9.10'
$ws.Range('E12').Value = '  +3.04%  '
$ws.Range('E13').Value = '  +1.55%  '
$ws.Range('D14').Value = "'" + '15.32'
$ws.Range('E14').Value = '  -2.44%  '
$ws.Range('D15').Value = '2.611.67'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = "'" + '0.856'
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('D17').Value = '2.268.39'
$ws.Range('E17').Value = '  -1.03%  '
$ws.Range('D18').Value = '43.086.53'
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('D19').Value = "'" + '0.0000107'
$ws.Range('E19').Value = '  -2.69%  '
$ws.Range('D20').Value = "'" + '6.70'
$ws.Range('E20').Value = '  +2.28%  '
$ws.Range('D21').Value = "'" + '71.72'
$ws.Range('E21').Value = '  -1.03%  '
$ws.Range('D22').Value = "'" + '2.42'
$ws.Range('E22').Value = '  -2.63%  '
$ws.Range('D23').Value = "'" + '233.47'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').Value = "'" + '9.37'
$ws.Range('E24').Value = '  -3.97%  '
$ws.Range('D25').Value = "'" + '2.83'
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  +1.88%  '
$ws.Range('D27').Value = "'" + '11.25'
$ws.Range('E27').Value = '  -3.46%  '
$ws.Range('D28').Value = "'" + '40.74'
$ws.Range('E28').Value = '  -3.28%  '
$ws.Range('E29').Value = '  -1.68%  '
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('D31').Value = "'" + '173.13'
$ws.Range('E31').Value = '  -2.23%  '
$ws.Range('D32').Value = "'" + '21.33'
$ws.Range('E32').Value = '  -1.30%  '
$ws.Range('D33').Value = "'" + '0.0895'
$ws.Range('E33').Value = '  -3.80%  '
$ws.Range('D34').Value = "'" + '5.60'
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('E35').Value = '  +2.49%  '
$ws.Range('D36').Value = "'" + '0.0368'
$ws.Range('E36').Value = '  +1.90%  '
$ws.Range('D37').Value = "'" + '4.59'
$ws.Range('E37').Value = '  -3.72%  '
$ws.Range('D38').Value = "'" + '3.90'
$ws.Range('E38').Value = '  +2.90%  '
$ws.Range('D39').Value = "'" + '0.104'
$ws.Range('E39').Value = '  -3.96%  '
$ws.Range('E40').Value = '  +7.14%  '
$ws.Range('D41').Value = "'" + '14.24'
$ws.Range('E41').Value = '  +1.74%  '
$ws.Range('D42').Value = "'" + '75.37'
$ws.Range('E42').Value = '  +5.23%  '
$ws.Range('D43').Value = "'" + '0.235'
$ws.Range('E43').Value = '  -3.97%  '
$ws.Range('D44').Value = "'" + '6.05'
$ws.Range('E44').Value = '  -1.29%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  -4.47%  '
$ws.Range('D47').Value = "'" + '8.48'
$ws.Range('E47').Value = '  -3.72%  '
$ws.Range('E48').Value = '  +1.48%  '
$ws.Range('D49').Value = "'" + '0.0993'
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('D50').Value = "'" + '100.48'
$ws.Range('E50').Value = '  -2.37%  '
$ws.Range('D51').Value = "'" + '0.593'
$ws.Range('E51').Value = '  +7.41%  '
